$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing A/B/C -> B/C/D
$ws.Columns.Item(1).Insert()

# Header row
$ws.Range("A1").Value = "Topic"
$ws.Range("B1").Value = "Words"
$ws.Range("C1").Value = "Frequency"
$ws.Range("D1").Value = "Name"

# Match the header formatting already applied to the other header cells
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$data = @(
    @(1,  "konut, deprem, temel, at, ev",                     98, "konut, deprem, temel, at, ev"),
    @(0,  "muhteşem, buluş, millet, hatay, erzurum",           90, "buluş, millet, muhteşem, erzurum, hatay"),
    @(2,  "rahmet, atatürk, şehit, gazi, teşekkür",            54, "rahmet, atatürk, şehit, dile, allah"),
    @(3,  "türkiye, yüzyıl, ülke, mayıs, millet",               52, "türkiye, yüzyıl, ülke, mayıs, millet"),
    @(5,  "teşekkür, kardeş, muhteşem, güzel, istanbul",       41, "teşekkür, kardeş, muhteşem, güzel, istanbul"),
    @(4,  "seçim, sandık, mayıs, millet, oy",                  35, "seçim, sandık, millet, mayıs, oy"),
    @(8,  "milyar, yatırım, lira, kamu, yıl",                  28, "milyar, yatırım, lira, kamu, dolar"),
    @(9,  "nükleer, santral, enerji, üretim, ülke",             23, "nükleer, santral, enerji, üretim, ülke"),
    @(6,  "canlı, yayın, ortak, tv, bölüm",                    16, "canlı, yayın, tv, basın, açıkla"),
    @(10, "emekli, maaş, hanım, aile, düşük",                  14, "emekli, maaş, hanım, aile, düşük"),
    @(11, "uzay, altay, milli, uçak, tank",                    12, "uzay, altay, milli, tank, uçak"),
    @(7,  "togg, uğurla, yol, mersin, yalnız",                 11, "gel, sev, haydi, vatan, bura")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
